$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 2.940931333333333
$ws.Cells.Item(2, 8).Value = 8.822794
$ws.Cells.Item(2, 9).Value = 0.03454096854573427
$ws.Cells.Item(2, 10).Value = 0.03454096854573427
$ws.Cells.Item(2, 11).Value = 1.0
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.065175
$ws.Cells.Item(2, 14).Value = 0.195525
$ws.Cells.Item(2, 15).Value = 0.009404016458916581
$ws.Cells.Item(2, 16).Value = 0.009404016458916581
$ws.Cells.Item(2, 17).Value = 0.19167519965
$ws.Cells.Item(2, 18).Value = 1.72507679685
$ws.Cells.Item(2, 19).Value = 0.000324823836711005
$ws.Cells.Item(2, 20).Value = 0.000324823836711005
$ws.Cells.Item(3, 7).Value = 2.940931333333333
$ws.Cells.Item(3, 8).Value = 8.822794
$ws.Cells.Item(3, 9).Value = 0.03454096854573427
$ws.Cells.Item(3, 10).Value = 0.03454096854573427
$ws.Cells.Item(3, 13).Value = 6.718514333333332
$ws.Cells.Item(3, 14).Value = 20.155543
$ws.Cells.Item(3, 15).Value = 0.969405744075698
$ws.Cells.Item(3, 16).Value = 0.969405744075698
$ws.Cells.Item(3, 17).Value = 19.75868931634911
$ws.Cells.Item(3, 18).Value = 177.828203847142
$ws.Cells.Item(3, 19).Value = 0.03348421331417281
$ws.Cells.Item(3, 20).Value = 0.03348421331417281
$ws.Cells.Item(4, 7).Value = 2.940931333333333
$ws.Cells.Item(4, 8).Value = 8.822794
$ws.Cells.Item(4, 9).Value = 0.03454096854573427
$ws.Cells.Item(4, 10).Value = 0.03454096854573427
$ws.Cells.Item(4, 11).Value = 1.0
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.14686
$ws.Cells.Item(4, 14).Value = 0.44058
$ws.Cells.Item(4, 15).Value = 0.02119023946538534
$ws.Cells.Item(4, 16).Value = 0.02119023946538533
$ws.Cells.Item(4, 17).Value = 0.4319051756133334
$ws.Cells.Item(4, 18).Value = 3.88714658052
$ws.Cells.Item(4, 19).Value = 0.0007319313948504519
$ws.Cells.Item(4, 20).Value = 0.0007319313948504517
$ws.Cells.Item(5, 9).Value = 0.8984069237831173
$ws.Cells.Item(5, 10).Value = 0.8984069237831174
$ws.Cells.Item(5, 11).Value = 1.0
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.065175
$ws.Cells.Item(5, 14).Value = 0.195525
$ws.Cells.Item(5, 15).Value = 0.009404016458916581
$ws.Cells.Item(5, 16).Value = 0.009404016458916581
$ws.Cells.Item(5, 17).Value = 4.985451587875
$ws.Cells.Item(5, 18).Value = 44.869064290875
$ws.Cells.Item(5, 19).Value = 0.00844863349806105
$ws.Cells.Item(5, 20).Value = 0.00844863349806105
$ws.Cells.Item(6, 9).Value = 0.8984069237831173
$ws.Cells.Item(6, 10).Value = 0.8984069237831174
$ws.Cells.Item(6, 13).Value = 6.718514333333332
$ws.Cells.Item(6, 14).Value = 20.155543
$ws.Cells.Item(6, 15).Value = 0.969405744075698
$ws.Cells.Item(6, 16).Value = 0.969405744075698
$ws.Cells.Item(6, 17).Value = 513.9214108366339
$ws.Cells.Item(6, 18).Value = 4625.292697529705
$ws.Cells.Item(6, 19).Value = 0.8709208324327317
$ws.Cells.Item(6, 20).Value = 0.8709208324327318
$ws.Cells.Item(7, 9).Value = 0.8984069237831173
$ws.Cells.Item(7, 10).Value = 0.8984069237831174
$ws.Cells.Item(7, 11).Value = 1.0
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.14686
$ws.Cells.Item(7, 14).Value = 0.44058
$ws.Cells.Item(7, 15).Value = 0.02119023946538534
$ws.Cells.Item(7, 16).Value = 0.02119023946538533
$ws.Cells.Item(7, 17).Value = 11.23380775136667
$ws.Cells.Item(7, 18).Value = 101.1042697623
$ws.Cells.Item(7, 19).Value = 0.01903745785232445
$ws.Cells.Item(7, 20).Value = 0.01903745785232445
$ws.Cells.Item(8, 7).Value = 1.002166333333333
$ws.Cells.Item(8, 8).Value = 3.006499
$ws.Cells.Item(8, 9).Value = 0.01177035159063915
$ws.Cells.Item(8, 10).Value = 0.01177035159063915
$ws.Cells.Item(8, 11).Value = 1.0
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.065175
$ws.Cells.Item(8, 14).Value = 0.195525
$ws.Cells.Item(8, 15).Value = 0.009404016458916581
$ws.Cells.Item(8, 16).Value = 0.009404016458916581
$ws.Cells.Item(8, 17).Value = 0.06531619077499999
$ws.Cells.Item(8, 18).Value = 0.587845716975
$ws.Cells.Item(8, 19).Value = 0.0001106885800856055
$ws.Cells.Item(8, 20).Value = 0.0001106885800856055
$ws.Cells.Item(9, 7).Value = 1.002166333333333
$ws.Cells.Item(9, 8).Value = 3.006499
$ws.Cells.Item(9, 9).Value = 0.01177035159063915
$ws.Cells.Item(9, 10).Value = 0.01177035159063915
$ws.Cells.Item(9, 13).Value = 6.718514333333332
$ws.Cells.Item(9, 14).Value = 20.155543
$ws.Cells.Item(9, 15).Value = 0.969405744075698
$ws.Cells.Item(9, 16).Value = 0.969405744075698
$ws.Cells.Item(9, 17).Value = 6.733068874884109
$ws.Cells.Item(9, 18).Value = 60.59761987395699
$ws.Cells.Item(9, 19).Value = 0.01141024644175612
$ws.Cells.Item(9, 20).Value = 0.01141024644175612
$ws.Cells.Item(10, 7).Value = 1.002166333333333
$ws.Cells.Item(10, 8).Value = 3.006499
$ws.Cells.Item(10, 9).Value = 0.01177035159063915
$ws.Cells.Item(10, 10).Value = 0.01177035159063915
$ws.Cells.Item(10, 11).Value = 1.0
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.14686
$ws.Cells.Item(10, 14).Value = 0.44058
$ws.Cells.Item(10, 15).Value = 0.02119023946538534
$ws.Cells.Item(10, 16).Value = 0.02119023946538533
$ws.Cells.Item(10, 17).Value = 0.1471781477133333
$ws.Cells.Item(10, 18).Value = 1.32460332942
$ws.Cells.Item(10, 19).Value = 0.0002494165687974228
$ws.Cells.Item(10, 20).Value = 0.0002494165687974227
$ws.Cells.Item(11, 7).Value = 2.356521666666667
$ws.Cells.Item(11, 8).Value = 7.069565
$ws.Cells.Item(11, 9).Value = 0.02767713065691252
$ws.Cells.Item(11, 10).Value = 0.02767713065691253
$ws.Cells.Item(11, 11).Value = 1.0
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.065175
$ws.Cells.Item(11, 14).Value = 0.195525
$ws.Cells.Item(11, 15).Value = 0.009404016458916581
$ws.Cells.Item(11, 16).Value = 0.009404016458916581
$ws.Cells.Item(11, 17).Value = 0.153586299625
$ws.Cells.Item(11, 18).Value = 1.382276696625
$ws.Cells.Item(11, 19).Value = 0.0002602761922331901
$ws.Cells.Item(11, 20).Value = 0.0002602761922331901
$ws.Cells.Item(12, 7).Value = 2.356521666666667
$ws.Cells.Item(12, 8).Value = 7.069565
$ws.Cells.Item(12, 9).Value = 0.02767713065691252
$ws.Cells.Item(12, 10).Value = 0.02767713065691253
$ws.Cells.Item(12, 13).Value = 6.718514333333332
$ws.Cells.Item(12, 14).Value = 20.155543
$ws.Cells.Item(12, 15).Value = 0.969405744075698
$ws.Cells.Item(12, 16).Value = 0.969405744075698
$ws.Cells.Item(12, 17).Value = 15.83232459431055
$ws.Cells.Item(12, 18).Value = 142.490921348795
$ws.Cells.Item(12, 19).Value = 0.0268303694383446
$ws.Cells.Item(12, 20).Value = 0.0268303694383446
$ws.Cells.Item(13, 7).Value = 2.356521666666667
$ws.Cells.Item(13, 8).Value = 7.069565
$ws.Cells.Item(13, 9).Value = 0.02767713065691252
$ws.Cells.Item(13, 10).Value = 0.02767713065691253
$ws.Cells.Item(13, 11).Value = 1.0
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.14686
$ws.Cells.Item(13, 14).Value = 0.44058
$ws.Cells.Item(13, 15).Value = 0.02119023946538534
$ws.Cells.Item(13, 16).Value = 0.02119023946538533
$ws.Cells.Item(13, 17).Value = 0.3460787719666667
$ws.Cells.Item(13, 18).Value = 3.1147089477
$ws.Cells.Item(13, 19).Value = 0.0005864850263347342
$ws.Cells.Item(13, 20).Value = 0.0005864850263347341
$ws.Cells.Item(14, 5).Value = 2.0
$ws.Cells.Item(14, 6).Value = 0.6666666666666666
$ws.Cells.Item(14, 7).Value = 0.5156633333333334
$ws.Cells.Item(14, 8).Value = 1.54699
$ws.Cells.Item(14, 9).Value = 0.006056418514425867
$ws.Cells.Item(14, 10).Value = 0.006056418514425868
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.065175
$ws.Cells.Item(14, 14).Value = 0.195525
$ws.Cells.Item(14, 15).Value = 0.009404016458916581
$ws.Cells.Item(14, 16).Value = 0.009404016458916581
$ws.Cells.Item(14, 17).Value = 0.03360835775
$ws.Cells.Item(14, 18).Value = 0.30247521975
$ws.Cells.Item(14, 19).Value = 0.00005695465939174797
$ws.Cells.Item(14, 20).Value = 0.00005695465939174797
$ws.Cells.Item(15, 5).Value = 2.0
$ws.Cells.Item(15, 6).Value = 0.6666666666666666
$ws.Cells.Item(15, 7).Value = 0.5156633333333334
$ws.Cells.Item(15, 8).Value = 1.54699
$ws.Cells.Item(15, 9).Value = 0.006056418514425867
$ws.Cells.Item(15, 10).Value = 0.006056418514425868
$ws.Cells.Item(15, 13).Value = 6.718514333333332
$ws.Cells.Item(15, 14).Value = 20.155543
$ws.Cells.Item(15, 15).Value = 0.969405744075698
$ws.Cells.Item(15, 16).Value = 0.969405744075698
$ws.Cells.Item(15, 17).Value = 3.464491496174444
$ws.Cells.Item(15, 18).Value = 31.18042346557
$ws.Cells.Item(15, 19).Value = 0.005871126896410841
$ws.Cells.Item(15, 20).Value = 0.005871126896410842
$ws.Cells.Item(16, 5).Value = 2.0
$ws.Cells.Item(16, 6).Value = 0.6666666666666666
$ws.Cells.Item(16, 7).Value = 0.5156633333333334
$ws.Cells.Item(16, 8).Value = 1.54699
$ws.Cells.Item(16, 9).Value = 0.006056418514425867
$ws.Cells.Item(16, 10).Value = 0.006056418514425868
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.14686
$ws.Cells.Item(16, 14).Value = 0.44058
$ws.Cells.Item(16, 15).Value = 0.02119023946538534
$ws.Cells.Item(16, 16).Value = 0.02119023946538533
$ws.Cells.Item(16, 17).Value = 0.07573031713333335
$ws.Cells.Item(16, 18).Value = 0.6815728542000001
$ws.Cells.Item(16, 19).Value = 0.0001283369586232775
$ws.Cells.Item(16, 20).Value = 0.0001283369586232774
$ws.Cells.Item(17, 5).Value = 3.0
$ws.Cells.Item(17, 6).Value = 1.0
$ws.Cells.Item(17, 7).Value = 1.834685
$ws.Cells.Item(17, 8).Value = 5.504055
$ws.Cells.Item(17, 9).Value = 0.02154820690917088
$ws.Cells.Item(17, 10).Value = 0.02154820690917089
$ws.Cells.Item(17, 11).Value = 1.0
$ws.Cells.Item(17, 12).Value = 0.3333333333333333
$ws.Cells.Item(17, 13).Value = 0.065175
$ws.Cells.Item(17, 14).Value = 0.195525
$ws.Cells.Item(17, 15).Value = 0.009404016458916581
$ws.Cells.Item(17, 16).Value = 0.009404016458916581
$ws.Cells.Item(17, 17).Value = 0.119575594875
$ws.Cells.Item(17, 18).Value = 1.076180353875
$ws.Cells.Item(17, 19).Value = 0.000202639692433983
$ws.Cells.Item(17, 20).Value = 0.000202639692433983
$ws.Cells.Item(18, 5).Value = 3.0
$ws.Cells.Item(18, 6).Value = 1.0
$ws.Cells.Item(18, 7).Value = 1.834685
$ws.Cells.Item(18, 8).Value = 5.504055
$ws.Cells.Item(18, 9).Value = 0.02154820690917088
$ws.Cells.Item(18, 10).Value = 0.02154820690917089
$ws.Cells.Item(18, 13).Value = 6.718514333333332
$ws.Cells.Item(18, 14).Value = 20.155543
$ws.Cells.Item(18, 15).Value = 0.969405744075698
$ws.Cells.Item(18, 16).Value = 0.969405744075698
$ws.Cells.Item(18, 17).Value = 12.32635746965166
$ws.Cells.Item(18, 18).Value = 110.937217226865
$ws.Cells.Item(18, 19).Value = 0.0208889555522819
$ws.Cells.Item(18, 20).Value = 0.0208889555522819
$ws.Cells.Item(19, 5).Value = 3.0
$ws.Cells.Item(19, 6).Value = 1.0
$ws.Cells.Item(19, 7).Value = 1.834685
$ws.Cells.Item(19, 8).Value = 5.504055
$ws.Cells.Item(19, 9).Value = 0.02154820690917088
$ws.Cells.Item(19, 10).Value = 0.02154820690917089
$ws.Cells.Item(19, 11).Value = 1.0
$ws.Cells.Item(19, 12).Value = 0.3333333333333333
$ws.Cells.Item(19, 13).Value = 0.14686
$ws.Cells.Item(19, 14).Value = 0.44058
$ws.Cells.Item(19, 15).Value = 0.02119023946538534
$ws.Cells.Item(19, 16).Value = 0.02119023946538533
$ws.Cells.Item(19, 17).Value = 0.2694418391
$ws.Cells.Item(19, 18).Value = 2.4249765519
$ws.Cells.Item(19, 19).Value = 0.0004566116644550019
$ws.Cells.Item(19, 20).Value = 0.0004566116644550019
